$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear contents of the cells that were removed in the diff
$ws.Range("C8").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("E14").ClearContents()

# Update the view: scroll position (topLeftCell A1 -> A5) and selection
# (activeCell I6/sqref I6 -> activeCell E12/sqref E12:E14)
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E12:E14").Select()
